$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44883
$ws.Range("J3").Value = 290
$ws.Range("K3").Value = 1400
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1434
$ws.Range("P3").Value = 1434

# Row 4
$ws.Range("D4").Value = 44210
$ws.Range("J4").Value = 1450
$ws.Range("K4").Value = 1600
$ws.Range("L4").Value = 1700
$ws.Range("M4").Value = 1650
$ws.Range("P4").Value = 1650

# Row 6
$ws.Range("D6").Value = 44537
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 1400
$ws.Range("M6").Value = 1350
$ws.Range("P6").Value = 1350

# Row 7
$ws.Range("D7").Value = 44638
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2800
$ws.Range("M7").Value = 2650
$ws.Range("P7").Value = 2650
